$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values that look numeric stay as text, matching original inline string formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "36.574.54"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "2.023.22"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "234.81"
$ws.Range("E5").Value = "  -9.56%  "
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "54.81"
$ws.Range("E8").Value = "  -3.03%  "
$ws.Range("E9").Value = "  -3.08%  "
$ws.Range("D10").Value = "58.11"
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").Value = "0.0748"
$ws.Range("E11").Value = "  -2.97%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "2.319.87"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "14.17"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "20.12"
$ws.Range("E15").Value = "  -7.09%  "
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "5.10"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "2.026.81"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "36.552.49"
$ws.Range("E19").Value = "  -1.95%  "
$ws.Range("D20").Value = "67.76"
$ws.Range("E20").Value = "  -3.18%  "
$ws.Range("D21").Value = "0.0₃0797"
$ws.Range("E21").Value = "  -4.45%  "
$ws.Range("E22").Value = "  +5.04%  "
$ws.Range("D23").Value = "220.35"
$ws.Range("E23").Value = "  -5.09%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  +1.25%  "
$ws.Range("D26").Value = "2.39"
$ws.Range("E26").Value = "  -7.43%  "
$ws.Range("D27").Value = "163.75"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "8.62"
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").Value = "18.94"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "4.36"
$ws.Range("E33").Value = "  -4.85%  "
$ws.Range("E34").Value = "  -5.48%  "
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("E36").Value = "  -4.80%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  -2.76%  "
$ws.Range("D40").Value = "5.70"
$ws.Range("E40").Value = "  +4.68%  "
$ws.Range("D41").Value = "2.97"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "4.33"
$ws.Range("E42").Value = "  +46.45%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.456.26"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0931"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("E46").Value = "  -6.14%  "
$ws.Range("D47").Value = "90.04"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").Value = "15.29"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").Value = "2.87"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -1.80%  "
